$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1405.6666
$ws.Range("J19").Value = 378.5
$ws.Range("L19").Value = 378.5
$ws.Range("N19").Value = -728.5
$ws.Range("H64").Value = 4025
$ws.Range("H67").Value = 4025
$ws.Range("H76").Value = 3000
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 1000
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 1000
$ws.Range("M76").Value = -4685
$ws.Range("N76").Value = -1630
$ws.Range("H79").Value = 3000
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 1000
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 1000
$ws.Range("M79").Value = -3908
$ws.Range("N79").Value = -3184
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H92").Value = 2381.6
$ws.Range("I92").Value = 301.33334
$ws.Range("K92").Value = 301.33334
$ws.Range("M92").Value = 946.66666
$ws.Range("H98").Value = 13326.223
$ws.Range("I98").Value = 11287.2
$ws.Range("J98").Value = 15875
$ws.Range("K98").Value = 11287.2
$ws.Range("L98").Value = 15875
$ws.Range("M98").Value = -9789.200000000001
$ws.Range("N98").Value = -18871
$ws.Range("H122").Value = 13326.223
$ws.Range("I122").Value = 11287.2
$ws.Range("J122").Value = 15875
$ws.Range("K122").Value = 33861.60000000001
$ws.Range("L122").Value = 47625
$ws.Range("M122").Value = -31411.60000000001
$ws.Range("N122").Value = -52525
$ws.Range("H138").Value = 7695968
$ws.Range("J138").Value = 3948.125
$ws.Range("L138").Value = 11844.375
$ws.Range("N138").Value = -22124.375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6267.8
$ws.Range("I63").Value = 6267.8
$ws.Range("K63").Value = 6267.8
$ws.Range("M63").Value = -5581.8
$ws.Range("H66").Value = 6267.8
$ws.Range("I66").Value = 6267.8
$ws.Range("K66").Value = 31339
$ws.Range("M66").Value = -27907
$ws.Range("H88").Value = 966.6667
$ws.Range("J88").Value = 900
$ws.Range("L88").Value = 900
$ws.Range("N88").Value = -1712
$ws.Range("H91").Value = 966.6667
$ws.Range("J91").Value = 900
$ws.Range("L91").Value = 900
$ws.Range("N91").Value = -3708
$ws.Range("H128").Value = 9999
$ws.Range("J128").Value = 9999
$ws.Range("L128").Value = 9999
$ws.Range("N128").Value = -19959
$ws.Range("H132").Value = 6464.3335
$ws.Range("I132").Value = 4739.857
$ws.Range("K132").Value = 14219.571
$ws.Range("M132").Value = -11689.571

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1415.875
$ws.Range("I86").Value = 1415.875
$ws.Range("K86").Value = 1415.875
$ws.Range("M86").Value = -292.875
$ws.Range("H89").Value = 1415.875
$ws.Range("I89").Value = 1415.875
$ws.Range("K89").Value = 7079.375
$ws.Range("M89").Value = -1463.375
$ws.Range("H94").Value = 817
$ws.Range("I94").Value = 817
$ws.Range("K94").Value = 817
$ws.Range("M94").Value = -366
$ws.Range("H134").Value = 6435.75
$ws.Range("I134").Value = 4581.3335
$ws.Range("K134").Value = 13744.0005
$ws.Range("M134").Value = -11209.0005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2583
$ws.Range("I62").Value = 2583
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2583
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1959
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2583
$ws.Range("I65").Value = 2583
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 12915
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -9795
$ws.Range("N65").ClearContents()
$ws.Range("H105").Value = 706.4
$ws.Range("I105").Value = 799.5
$ws.Range("K105").Value = 799.5
$ws.Range("M105").Value = 947.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 71.666664
$ws.Range("I33").Value = 82.5
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 495
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -212
$ws.Range("N33").Value = -866
$ws.Range("H68").Value = 261
$ws.Range("I68").Value = 22
$ws.Range("J68").Value = 500
$ws.Range("K68").Value = 66
$ws.Range("L68").Value = 1500
$ws.Range("M68").Value = 745
$ws.Range("N68").Value = -3122
$ws.Range("H71").Value = 261
$ws.Range("I71").Value = 22
$ws.Range("J71").Value = 500
$ws.Range("K71").Value = 198
$ws.Range("L71").Value = 4500
$ws.Range("M71").Value = 3858
$ws.Range("N71").Value = -12612
$ws.Range("H131").Value = 923.25
$ws.Range("J131").Value = 533
$ws.Range("L131").Value = 1599
$ws.Range("N131").Value = -11679

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3749
$ws.Range("I70").Value = 3749
$ws.Range("K70").Value = 3749
$ws.Range("M70").Value = -3479
$ws.Range("H73").Value = 3749
$ws.Range("I73").Value = 3749
$ws.Range("K73").Value = 3749
$ws.Range("M73").Value = -2813
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H113").Value = 3132.375
$ws.Range("I113").Value = 3132.375
$ws.Range("K113").Value = 3132.375
$ws.Range("M113").Value = -962.375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 600.5
$ws.Range("J16").Value = 1002
$ws.Range("L16").Value = 1002
$ws.Range("N16").Value = -1342
$ws.Range("H46").Value = 7467
$ws.Range("I46").Value = 7450.5
$ws.Range("K46").Value = 7450.5
$ws.Range("M46").Value = -7262.5
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 999.5
$ws.Range("I81").Value = 999.5
$ws.Range("K81").Value = 1999
$ws.Range("M81").Value = -938
$ws.Range("H84").Value = 999.5
$ws.Range("I84").Value = 999.5
$ws.Range("K84").Value = 9995
$ws.Range("M84").Value = -4691
$ws.Range("H136").Value = 3521.4285
$ws.Range("I136").Value = 3982
$ws.Range("K136").Value = 11946
$ws.Range("M136").Value = -9396
